$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.315.42'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '3.398.45'
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.03'
$ws.Range("E5").Value = '  +1.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.71'
$ws.Range("E6").Value = '  +2.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '3.396.51'
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.546'
$ws.Range("E9").Value = '  +2.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.39'
$ws.Range("E10").Value = '  -0.99%  '
$ws.Range("E11").Value = '  +3.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.432'
$ws.Range("E12").Value = '  -0.88%  '
$ws.Range("D13").Value = '3.985.55'
$ws.Range("E13").Value = '  +1.69%  '
$ws.Range("E14").Value = '  -3.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000193'
$ws.Range("E15").Value = '  +7.37%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.23'
$ws.Range("E16").Value = '  +1.26%  '
$ws.Range("D17").Value = '63.371.03'
$ws.Range("E17").Value = '  +1.05%  '
$ws.Range("D18").Value = '3.367.58'
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.25'
$ws.Range("E19").Value = '  -1.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.07'
$ws.Range("E20").Value = '  +1.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '379.43'
$ws.Range("E21").Value = '  -1.51%  '
$ws.Range("E22").Value = '  -3.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.46'
$ws.Range("E24").Value = '  +1.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.529'
$ws.Range("E25").Value = '  -1.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000120'
$ws.Range("E26").Value = '  +26.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.43'
$ws.Range("E27").Value = '  +6.21%  '
$ws.Range("E28").Value = '  -0.10%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("E30").Value = '  +8.40%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.36'
$ws.Range("E31").Value = '  +3.87%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.00'
$ws.Range("E32").Value = '  +0.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.17'
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.37'
$ws.Range("E34").Value = '  -2.91%  '
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.81'
$ws.Range("E36").Value = '  +1.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '159.69'
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.44'
$ws.Range("E38").Value = '  -1.46%  '
$ws.Range("D39").Value = '2.945.25'
$ws.Range("E39").Value = '  +4.46%  '
$ws.Range("E40").Value = '  +0.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.83'
$ws.Range("E41").Value = '  -3.30%  '
$ws.Range("E42").Value = '  +1.13%  '
$ws.Range("E43").Value = '  +1.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.75'
$ws.Range("E44").Value = '  +3.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.762'
$ws.Range("E45").Value = '  +2.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.33'
$ws.Range("E46").Value = '  +1.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.23'
$ws.Range("E47").Value = '  +6.23%  '
$ws.Range("E48").Value = '  +3.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.20'
$ws.Range("E49").Value = '  +22.77%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.840'
$ws.Range("E50").Value = '  +5.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.36'
$ws.Range("E51").Value = '  +1.16%  '
